$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 28; existing rows 28..84 shift down to 29..85.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with the new data record.
$ws.Cells.Item(28, 1).Value = 11
$ws.Cells.Item(28, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(28, 3).Value = "Bíobío"
$ws.Cells.Item(28, 4).Value = (Get-Date -Year 2021 -Month 12 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(28, 5).Value = 8
$ws.Cells.Item(28, 6).Value = "Fruta"
$ws.Cells.Item(28, 7).Value = 100108
$ws.Cells.Item(28, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(28, 9).Value = 100108002
$ws.Cells.Item(28, 10).Value = "Mango"
$ws.Cells.Item(28, 11).Value = "Sin especificar"
$ws.Cells.Item(28, 12).Value = "Primera"
$ws.Cells.Item(28, 13).Value = 200
$ws.Cells.Item(28, 14).Value = 6000
$ws.Cells.Item(28, 15).Value = 6500
$ws.Cells.Item(28, 16).Value = 6250
$ws.Cells.Item(28, 17).Value = "`$/bandeja 4 kilos"
$ws.Cells.Item(28, 18).Value = "Perú"
$ws.Cells.Item(28, 19).Value = 1562
$ws.Cells.Item(28, 20).Value = 4
